$wb = $excel.ActiveWorkbook

# Sheet2 gets brand-new content: a "sheet1"/"zhu" header row and a data row,
# replacing the old "sheet2"/"sheet2" header + lone B2 value.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Value = "sheet1"
$ws2.Range("B1").Value = "zhu"
$ws2.Range("A2").Value = 3.3
$ws2.Range("B2").Value = 100

# Sheet3 loses its second row (B2) so it only keeps the header row.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B2").ClearContents()
$null = $ws3.Range("B6").Select()

# Sheet2 becomes the active/selected sheet (was Sheet3), with B2 selected.
$null = $ws2.Activate()
$null = $ws2.Range("B2").Select()
